# Fix Training Data Issue
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") contained the malformed string "6-13-2013-14" for
# every data row; correct it to the properly formatted date "2014-06-13".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-13-2013-14"
$newValue = "2014-06-13"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF" + $row)
    if ($cell.Text -eq $oldValue) {
        # Leading apostrophe forces Excel to keep the corrected date as
        # literal text instead of auto-converting it to a date serial.
        $cell.Value = "'" + $newValue
    }
}
